$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsA = @(
    45968,
    45968.01041666666,
    45968.02083333334,
    45968.03125,
    45968.04166666666,
    45968.05208333334,
    45968.0625,
    45968.07291666666,
    45968.08333333334,
    45968.09375,
    45968.10416666666,
    45968.11458333334,
    45968.125,
    45968.13541666666,
    45968.14583333334,
    45968.15625,
    45968.16666666666,
    45968.17708333334,
    45968.1875,
    45968.19791666666,
    45968.20833333334,
    45968.21875,
    45968.22916666666,
    45968.23958333334,
    45968.25,
    45968.26041666666,
    45968.27083333334,
    45968.28125,
    45968.29166666666,
    45968.30208333334,
    45968.3125,
    45968.32291666666,
    45968.33333333334,
    45968.34375,
    45968.35416666666,
    45968.36458333334,
    45968.375,
    45968.38541666666,
    45968.39583333334,
    45968.40625,
    45968.41666666666,
    45968.42708333334,
    45968.4375,
    45968.44791666666,
    45968.45833333334,
    45968.46875,
    45968.47916666666,
    45968.48958333334,
    45968.5,
    45968.51041666666,
    45968.52083333334,
    45968.53125,
    45968.54166666666,
    45968.55208333334,
    45968.5625,
    45968.57291666666,
    45968.58333333334,
    45968.59375,
    45968.60416666666,
    45968.61458333334,
    45968.625,
    45968.63541666666,
    45968.64583333334,
    45968.65625,
    45968.66666666666,
    45968.67708333334,
    45968.6875,
    45968.69791666666,
    45968.70833333334,
    45968.71875,
    45968.72916666666,
    45968.73958333334,
    45968.75,
    45968.76041666666,
    45968.77083333334,
    45968.78125,
    45968.79166666666,
    45968.80208333334,
    45968.8125,
    45968.82291666666,
    45968.83333333334,
    45968.84375,
    45968.85416666666,
    45968.86458333334,
    45968.875,
    45968.88541666666,
    45968.89583333334,
    45968.90625,
    45968.91666666666,
    45968.92708333334,
    45968.9375,
    45968.94791666666,
    45968.95833333334,
    45968.96875,
    45968.97916666666,
    45968.98958333334,
    45969,
    45969.01041666666,
    45969.02083333334,
    45969.03125,
    45969.04166666666,
    45969.05208333334,
    45969.0625,
    45969.07291666666,
    45969.08333333334,
    45969.09375,
    45969.10416666666,
    45969.11458333334,
    45969.125,
    45969.13541666666,
    45969.14583333334,
    45969.15625,
    45969.16666666666,
    45969.17708333334,
    45969.1875,
    45969.19791666666,
    45969.20833333334,
    45969.21875,
    45969.22916666666,
    45969.23958333334,
    45969.25,
    45969.26041666666,
    45969.27083333334,
    45969.28125,
    45969.29166666666,
    45969.30208333334,
    45969.3125,
    45969.32291666666,
    45969.33333333334,
    45969.34375,
    45969.35416666666,
    45969.36458333334,
    45969.375,
    45969.38541666666,
    45969.39583333334,
    45969.40625,
    45969.41666666666,
    45969.42708333334,
    45969.4375,
    45969.44791666666,
    45969.45833333334,
    45969.46875,
    45969.47916666666,
    45969.48958333334,
    45969.5,
    45969.51041666666,
    45969.52083333334,
    45969.53125,
    45969.54166666666,
    45969.55208333334,
    45969.5625,
    45969.57291666666,
    45969.58333333334,
    45969.59375,
    45969.60416666666,
    45969.61458333334,
    45969.625,
    45969.63541666666,
    45969.64583333334,
    45969.65625,
    45969.66666666666,
    45969.67708333334,
    45969.6875,
    45969.69791666666,
    45969.70833333334,
    45969.71875,
    45969.72916666666,
    45969.73958333334,
    45969.75,
    45969.76041666666,
    45969.77083333334,
    45969.78125,
    45969.79166666666,
    45969.80208333334,
    45969.8125,
    45969.82291666666,
    45969.83333333334,
    45969.84375,
    45969.85416666666,
    45969.86458333334,
    45969.875,
    45969.88541666666,
    45969.89583333334,
    45969.90625,
    45969.91666666666,
    45969.92708333334,
    45969.9375,
    45969.94791666666,
    45969.95833333334,
    45969.96875,
    45969.97916666666,
    45969.98958333334
)
$rowsB = @(
    5505,
    5509,
    5516,
    5520,
    5416,
    5435,
    5420,
    5490,
    5387,
    5366,
    5406,
    5451,
    5455,
    5395,
    5476,
    5500,
    5670,
    5713,
    5853,
    5953,
    6183,
    6304,
    6473,
    6648,
    6787,
    6938,
    7135,
    7196,
    7288,
    7300,
    7316,
    7375,
    7350,
    7364,
    7311,
    7302,
    7251,
    7226,
    7153,
    7101,
    6815,
    6915,
    6911,
    6893,
    6804,
    6842,
    6866,
    6881,
    6916,
    7050,
    7078,
    7198,
    7239,
    7256,
    7267,
    7267,
    7314,
    7361,
    7477,
    7562,
    7527,
    7524,
    7626,
    7697,
    7818,
    7991,
    8060,
    8112,
    8013,
    8046,
    8000,
    7933,
    7895,
    7868,
    7885,
    7908,
    7755,
    7700,
    7596,
    7507,
    7377,
    7287,
    7140,
    7026,
    6837,
    6726,
    6563,
    6466,
    6329,
    6224,
    6105,
    6077,
    5920,
    5835,
    5744,
    5666,
    5600,
    5569,
    5557,
    5475,
    5469,
    5449,
    5412,
    5365,
    5403,
    5364,
    5399,
    5391,
    5423,
    5402,
    5421,
    5401,
    5457,
    5475,
    5481,
    5526,
    5485,
    5534,
    5588,
    5636,
    5744,
    5821,
    5843,
    5915,
    5985,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0
)
$rowsD = @(
    "07.11.20251",
    "07.11.20252",
    "07.11.20253",
    "07.11.20254",
    "07.11.20255",
    "07.11.20256",
    "07.11.20257",
    "07.11.20258",
    "07.11.20259",
    "07.11.202510",
    "07.11.202511",
    "07.11.202512",
    "07.11.202513",
    "07.11.202514",
    "07.11.202515",
    "07.11.202516",
    "07.11.202517",
    "07.11.202518",
    "07.11.202519",
    "07.11.202520",
    "07.11.202521",
    "07.11.202522",
    "07.11.202523",
    "07.11.202524",
    "07.11.202525",
    "07.11.202526",
    "07.11.202527",
    "07.11.202528",
    "07.11.202529",
    "07.11.202530",
    "07.11.202531",
    "07.11.202532",
    "07.11.202533",
    "07.11.202534",
    "07.11.202535",
    "07.11.202536",
    "07.11.202537",
    "07.11.202538",
    "07.11.202539",
    "07.11.202540",
    "07.11.202541",
    "07.11.202542",
    "07.11.202543",
    "07.11.202544",
    "07.11.202545",
    "07.11.202546",
    "07.11.202547",
    "07.11.202548",
    "07.11.202549",
    "07.11.202550",
    "07.11.202551",
    "07.11.202552",
    "07.11.202553",
    "07.11.202554",
    "07.11.202555",
    "07.11.202556",
    "07.11.202557",
    "07.11.202558",
    "07.11.202559",
    "07.11.202560",
    "07.11.202561",
    "07.11.202562",
    "07.11.202563",
    "07.11.202564",
    "07.11.202565",
    "07.11.202566",
    "07.11.202567",
    "07.11.202568",
    "07.11.202569",
    "07.11.202570",
    "07.11.202571",
    "07.11.202572",
    "07.11.202573",
    "07.11.202574",
    "07.11.202575",
    "07.11.202576",
    "07.11.202577",
    "07.11.202578",
    "07.11.202579",
    "07.11.202580",
    "07.11.202581",
    "07.11.202582",
    "07.11.202583",
    "07.11.202584",
    "07.11.202585",
    "07.11.202586",
    "07.11.202587",
    "07.11.202588",
    "07.11.202589",
    "07.11.202590",
    "07.11.202591",
    "07.11.202592",
    "07.11.202593",
    "07.11.202594",
    "07.11.202595",
    "07.11.202596",
    "08.11.20251",
    "08.11.20252",
    "08.11.20253",
    "08.11.20254",
    "08.11.20255",
    "08.11.20256",
    "08.11.20257",
    "08.11.20258",
    "08.11.20259",
    "08.11.202510",
    "08.11.202511",
    "08.11.202512",
    "08.11.202513",
    "08.11.202514",
    "08.11.202515",
    "08.11.202516",
    "08.11.202517",
    "08.11.202518",
    "08.11.202519",
    "08.11.202520",
    "08.11.202521",
    "08.11.202522",
    "08.11.202523",
    "08.11.202524",
    "08.11.202525",
    "08.11.202526",
    "08.11.202527",
    "08.11.202528",
    "08.11.202529",
    "08.11.202530",
    "08.11.202531",
    "08.11.202532",
    "08.11.202533",
    "08.11.202534",
    "08.11.202535",
    "08.11.202536",
    "08.11.202537",
    "08.11.202538",
    "08.11.202539",
    "08.11.202540",
    "08.11.202541",
    "08.11.202542",
    "08.11.202543",
    "08.11.202544",
    "08.11.202545",
    "08.11.202546",
    "08.11.202547",
    "08.11.202548",
    "08.11.202549",
    "08.11.202550",
    "08.11.202551",
    "08.11.202552",
    "08.11.202553",
    "08.11.202554",
    "08.11.202555",
    "08.11.202556",
    "08.11.202557",
    "08.11.202558",
    "08.11.202559",
    "08.11.202560",
    "08.11.202561",
    "08.11.202562",
    "08.11.202563",
    "08.11.202564",
    "08.11.202565",
    "08.11.202566",
    "08.11.202567",
    "08.11.202568",
    "08.11.202569",
    "08.11.202570",
    "08.11.202571",
    "08.11.202572",
    "08.11.202573",
    "08.11.202574",
    "08.11.202575",
    "08.11.202576",
    "08.11.202577",
    "08.11.202578",
    "08.11.202579",
    "08.11.202580",
    "08.11.202581",
    "08.11.202582",
    "08.11.202583",
    "08.11.202584",
    "08.11.202585",
    "08.11.202586",
    "08.11.202587",
    "08.11.202588",
    "08.11.202589",
    "08.11.202590",
    "08.11.202591",
    "08.11.202592",
    "08.11.202593",
    "08.11.202594",
    "08.11.202595",
    "08.11.202596"
)

$startRow = 2
for ($i = 0; $i -lt $rowsA.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rowsA[$i]
    $ws.Cells.Item($r, 2).Value = $rowsB[$i]
    $ws.Cells.Item($r, 4).Value = $rowsD[$i]
}

Write-Output "Updated $($rowsA.Length) rows"